$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "326.86"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-1.20%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.04%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.263"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.83%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08349"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.55%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.426"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.05%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.938"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-5.81%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9695"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.65%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.541"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-2.72%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1133"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.56%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1888"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.18%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09693"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.36%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04604"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-3.84%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.05%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001294"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "3.21%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005779"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.28%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.403"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.80%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3344"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.13%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.610"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-15.28%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.29%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2578"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.34%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04153"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.16%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001234"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-5.37%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004428"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "1.02%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001302"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.94%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0002983"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-20.04%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02679"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-0.37%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05556"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-1.76%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007835"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.97%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007338"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-2.63%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002043"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.55%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007828"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-5.59%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3503"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006862"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-2.22%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.39%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003489"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-0.97%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003535"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "40.61%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002103"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.39%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002003"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.39%"
